# Top level menu, Restriction based data filtration, reports, new home page
# view and, issue fixes.
#
# On the "test" sheet:
#  - C2 (Interviewer for the first data row) gets the fuller credit line
#    (adds a new shared string and repoints C2 at it).
#  - Columns B (Interviewee) and C (Interviewer) get explicit widths so the
#    new, longer text is readable.
#  - The view no longer parks the viewport scrolled out to AA1/AF29; it now
#    opens with the selection at C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the Interviewer value for row 2 (adds a new shared string).
$ws.Range("C2").Value = "Little Thunder, Julie Pearson;Finchum, Tanya; Bishop, Alex"

# Give the Interviewee / Interviewer columns enough room for the text.
$ws.Columns.Item(2).ColumnWidth = 23.333333333333332
$ws.Columns.Item(3).ColumnWidth = 81.83333333333333

# Reset the view: select C3 (this also clears the old scrolled-out
# topLeftCell position).
[void]$ws.Range("C3").Select()
